$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "JMLBARANG" column (column F) is no longer needed; delete it entirely.
# This shifts every column to its right one position to the left, which
# matches the header/value layout of the updated sheet.
$ws.Columns.Item(6).Delete()

# Restore the selection to match what was active when the sheet was saved.
$ws.Range("H11").Select()
